$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J2: date-type code updated (leading-zero code, keep as text)
$ws.Range("J2").Value = "'004"

# Column N2: report date updated
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Raw financial figures for the new reporting period
$ws.Range("O2").Value = 2283089248.81
$ws.Range("P2").Value = 271549353.03
$ws.Range("Q2").Value = 103313792.42
$ws.Range("S2").Value = 320751738.77
$ws.Range("U2").Value = 415787055.17
$ws.Range("W2").Value = 1025124293.73
$ws.Range("X2").Value = 283548927.65
$ws.Range("Z2").Value = 117712.5
$ws.Range("AB2").Value = 1257964955.08
$ws.Range("AF2").Value = 162.501113076
$ws.Range("AG2").Value = 44.9007542856

# Ratio columns no longer reported for this period - clear them out
$ws.Range("R2").ClearContents()
$ws.Range("T2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
